# Updated bug fix list.
# Mark the last two "Manage Volunteers" issue rows (11 and 12) as Fixed,
# with the fixed date, matching the pattern used in rows 2-5, and move
# the active selection to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Profile / "Volunteer Type displays..." -> Fixed, with date
$ws.Range("C11").Value = "x"
$ws.Range("C11").Style = $ws.Range("C2").Style
$ws.Range("D11").Value = 42859
$ws.Range("D11").Style = $ws.Range("D2").Style

# Row 12: Profile / "Manage Children button..." -> Fixed, with date
$ws.Range("C12").Value = "x"
$ws.Range("C12").Style = $ws.Range("C2").Style
$ws.Range("D12").Value = 42859
$ws.Range("D12").Style = $ws.Range("D2").Style

# Update the active cell selection to C13 (as recorded when the author saved)
$ws.Range("C13").Select()
